$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A3: remove leading non-breaking space ---
$ws.Range("A3").Value = "Methodist Ladies ' College"

# --- Copy row formats for the new rows ---
$ws.Range("A3:E3").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)

$ws.Range("A6:E6").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A11:E11").PasteSpecial(-4122)

# --- Row 7: Calvert / Charles_Calvert,_5th_Baron_Baltimore ---
$ws.Range("B7").Value = "53715bd0a310958a4e6adeb3"
$ws.Range("C7").Value = " Charles_Calvert,_5th_Baron_Baltimore"
$r7c = $ws.Range("C7").Characters(2, 36)
$r7c.Font.Name = "Times New Roman"
$r7c.Font.Size = 10
$r7c.Font.Color = 16711680
$ws.Range("D7").Value = "Baron_Baltimore"
$ws.Range("E7").Value = "link"
$ws.Hyperlinks.Add($ws.Range("C7"), "http://en.wikipedia.org/wiki/Charles_Calvert,_5th_Baron_Baltimore", "", "", "Charles_Calvert,_5th_Baron_Baltimore")
$ws.Hyperlinks.Add($ws.Range("D7"), "http://en.wikipedia.org/wiki/Baron_Baltimore", "", "", "Baron_Baltimore")
$ws.Hyperlinks.Add($ws.Range("E7"), "http://en.wikipedia.org/wiki/?curid=12859", "", "", "link")

# --- Row 8: Maccabi Haifa / Maccabi_Haifa_B.C. / Maccabi_Haifa_F.C. ---
$ws.Range("A8").Value = "Maccabi Haifa"
$ws.Range("B8").Value = "53715bd0a310958a4e6adeb8"
$ws.Range("C8").Value = "Maccabi_Haifa_B.C. "
$r8c = $ws.Range("C8").Characters(19, 1)
$r8c.Font.Name = "Times New Roman"
$r8c.Font.Size = 10
$r8c.Font.ColorIndex = -4105
$ws.Range("D8").Value = "Maccabi_Haifa_F.C."
$ws.Range("E8").Value = "link"
$ws.Hyperlinks.Add($ws.Range("C8"), "http://en.wikipedia.org/wiki/Maccabi_Haifa_B.C.", "", "", "Maccabi_Haifa_B.C.")
$ws.Hyperlinks.Add($ws.Range("D8"), "http://en.wikipedia.org/wiki/Maccabi_Haifa_F.C.", "", "", "Maccabi_Haifa_F.C.")
$ws.Hyperlinks.Add($ws.Range("E8"), "http://en.wikipedia.org/wiki/?curid=12859", "", "", "link")

# --- Row 9: British / United_Kingdom / Presidencies_and_provinces_of... ---
$ws.Range("A9").Value = "British"
$ws.Range("B9").Value = "53715bd0a310958a4e6adec4"
$ws.Range("C9").Value = "United_Kingdom "
$r9c = $ws.Range("C9").Characters(15, 1)
$r9c.Font.Name = "Times New Roman"
$r9c.Font.Size = 10
$r9c.Font.ColorIndex = -4105
$ws.Range("D9").Value = "Presidencies_and_provinces_of "
$ws.Range("E9").Value = "link"
$ws.Hyperlinks.Add($ws.Range("C9"), "http://en.wikipedia.org/wiki/United_Kingdom", "", "", "United_Kingdom")
$ws.Hyperlinks.Add($ws.Range("D9"), "http://en.wikipedia.org/wiki/Presidencies_and_provinces_of_British_India", "", "", "Presidencies_and_provinces_of ")
$ws.Hyperlinks.Add($ws.Range("E9"), "http://en.wikipedia.org/wiki/?curid=12859", "", "", "link")

# --- Row 10: triads / Triad_(music) / Chord_(music) ---
$ws.Range("A10").Value = "triads"
$ws.Range("B10").Value = "53715bd0a310958a4e6aded9"
$ws.Range("C10").Value = "Triad_(music)"
$ws.Range("D10").Value = "Chord_(music)"
$ws.Range("E10").Value = "link"
$ws.Hyperlinks.Add($ws.Range("C10"), "http://en.wikipedia.org/wiki/Triad_(music)", "", "", "Triad_(music)")
$ws.Hyperlinks.Add($ws.Range("D10"), "http://en.wikipedia.org/wiki/Chord_(music)", "", "", "Chord_(music)")
$ws.Hyperlinks.Add($ws.Range("E10"), "http://en.wikipedia.org/wiki/?curid=12859", "", "", "link")

# --- Row 11: Miesbach / Miesbach_(district) / Miesbach ---
$ws.Range("A11").Value = "Miesbach"
$ws.Range("B11").Value = "53715bd0a310958a4e6adee9"
$ws.Range("C11").Value = "Miesbach_(district) "
$r11c = $ws.Range("C11").Characters(20, 1)
$r11c.Font.Name = "Times New Roman"
$r11c.Font.Size = 10
$r11c.Font.ColorIndex = -4105
$ws.Range("D11").Value = "Miesbach"
$ws.Range("E11").Value = "link "
$ws.Hyperlinks.Add($ws.Range("C11"), "http://en.wikipedia.org/wiki/Miesbach_(district)", "", "", "Miesbach_(district)")
$ws.Hyperlinks.Add($ws.Range("D11"), "http://en.wikipedia.org/wiki/Miesbach", "", "", "Miesbach")
$ws.Hyperlinks.Add($ws.Range("E11"), "http://en.wikipedia.org/wiki/?curid=12859", "", "", "link")

# --- Column C width ---
$ws.Columns("C").ColumnWidth = 31.0765306122449

# --- Selection ---
$ws.Range("B11").Select() | Out-Null

Write-Host "All changes applied"
